# Tuntikirjanpito.xlsx - add a new time-tracking entry for the signup/login
# refactor work, then roll the running total forward to include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 70: 2 hours spent on signup/login refactor work (project: api)
$ws.Range("B70").Value = 2
$ws.Range("C70").Value = "login ja signup for refaktoroitu, uusi TextInput komponentti, signup formista puuttuu vielä errorit, signup service"
$ws.Range("D70").Value = "api"

# Extend the hours-total formula (row 75) to include the new row
$ws.Range("B75").Formula = "=SUM(B2:B70)"

# Move the active selection to the newly entered description cell
$ws.Range("C70").Select() | Out-Null
